$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column C for existing rows 1-18
$ws.Cells.Item(1, 3).Value = "ИСКО"
$ws.Cells.Item(2, 3).Value = "ЕБС"
$ws.Cells.Item(3, 3).Value = "ИСКО"
$ws.Cells.Item(4, 3).Value = "ИСКО"
$ws.Cells.Item(5, 3).Value = "ЕБС"
$ws.Cells.Item(6, 3).Value = "ИСКО"
$ws.Cells.Item(7, 3).Value = "ИСКО"
$ws.Cells.Item(8, 3).Value = "ИСКО"
$ws.Cells.Item(9, 3).Value = "ИСКО"
$ws.Cells.Item(10, 3).Value = "ИСКО"
$ws.Cells.Item(11, 3).Value = "ЕБС"
$ws.Cells.Item(12, 3).Value = "ИСКО"
$ws.Cells.Item(13, 3).Value = "ЕБС"
$ws.Cells.Item(14, 3).Value = "ЕБС"
$ws.Cells.Item(15, 3).Value = "ИСКО"
$ws.Cells.Item(16, 3).Value = "ИСКО"
$ws.Cells.Item(17, 3).Value = "ИСКО"
$ws.Cells.Item(18, 3).Value = "ИСКО"

# Insert 3 new rows at position 19, pushing old row 19 down to row 22
$ws.Rows("19:21").Insert()

# Fill the newly inserted rows 19-21
$ws.Cells.Item(19, 1).Value = "Неудовлетворительное качество БО"
$ws.Cells.Item(19, 2).Value = 422
$ws.Cells.Item(19, 3).Value = "ИСКО"
$ws.Cells.Item(20, 1).Value = "Ошибка при взаимодействии с ФХ СМЭВ"
$ws.Cells.Item(20, 2).Value = 432
$ws.Cells.Item(20, 3).Value = "СМЭВ"
$ws.Cells.Item(21, 1).Value = "Ошибка валидации логина пароля ФХ СМЭВ"
$ws.Cells.Item(21, 2).Value = 433
$ws.Cells.Item(21, 3).Value = "ИГНОРИРОВАТЬ"

# Set column C for the row that shifted down to 22 (Внутренняя ошибка сервера / 500)
$ws.Cells.Item(22, 3).Value = "ЕСИА"

# Add new row 23
$ws.Cells.Item(23, 1).Value = "Проблемы взаимодействия со СМЭВ (агрегированная ошибка)"
$ws.Cells.Item(23, 2).Value = 505

# Update selection to match target (D11) and dimension will auto-update
$ws.Range("D11").Select()
